$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row: insert a "day" (date) column at C, shift startTime/endTime
#     headers over, and rename the old "Space" header (now holding room data)
#     to lowercase "space" in column F ---
$ws.Range("C1").Value = "day"
$ws.Range("D1").Value = "startTime"
$ws.Range("E1").Value = "endTime"
$ws.Range("F1").Value = "space"

# --- Row 2: existing booking (Huzaifa RAGHAV) gets its startTime split into a
#     real date (column C) + a plain start-time string (column D); endTime and
#     room shift over to E/F ---
$ws.Range("C2").Value = 43046
$ws.Range("C2").NumberFormat = "d-mmm-yy"
$ws.Range("D2").Value = "12:45PM"
$ws.Range("E2").Value = "1:15PM"
$ws.Range("F2").Value = "Green Screen Room 1"

# --- Row 3: new booking (Carl CASTUERAS), written by the booking loop ---
$ws.Range("A3").Value = "Carl CASTUERAS"
$ws.Range("B3").Value = "cac59128@gapps.uwcsea.edu.sg"
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:cac59128@gapps.uwcsea.edu.sg")
$ws.Range("C3").Value = 43049
$ws.Range("C3").NumberFormat = "d-mmm-yy"
$ws.Range("D3").Value = "12:45PM"
$ws.Range("E3").Value = "1:15PM"
$ws.Range("F3").Value = "Green Screen Room 2"

$ws.Range("F4").Select() | Out-Null
